$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Three new "Sliding Window" problems were solved and logged in rows
# 56-58 (SubArray with Product Less then K, Max Consecutive Ones III,
# Fruits Into Basket). The three previously-pending "Sliding Window" /
# "Moore's Voting" / "DNF Algorithm" rows that used to sit at rows
# 56-60 were pushed down to rows 61-65, with a couple of blank spacer
# rows (59, 60 and 66) left in place, matching the new layout.
# ---------------------------------------------------------------------

# Row 56: edit in place - now "Done" with a new problem name & full details
$ws.Cells.Item(56, 3).Value = "SubArray with Product Less then K"
$ws.Cells.Item(56, 5).Value = "Done"
$ws.Cells.Item(56, 6).Value = 45898
$ws.Range("F55").Copy() | Out-Null
$ws.Range("F56").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(56, 7).Value = "O(n)"
$ws.Cells.Item(56, 8).Value = "O(1)"
$ws.Cells.Item(56, 9).Value = "Sliding Window"

# Row 57: edit in place - now "Done" with a new problem name & full details
$ws.Cells.Item(57, 3).Value = "Max Consecutive Ones III"
$ws.Cells.Item(57, 4).Value = "Medium"
$ws.Cells.Item(57, 5).Value = "Done"
$ws.Cells.Item(57, 6).Value = 45898
$ws.Range("F55").Copy() | Out-Null
$ws.Range("F57").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(57, 7).Value = "O(n)"
$ws.Cells.Item(57, 8).Value = "O(1)"
$ws.Cells.Item(57, 9).Value = "Sliding Window"

# Row 58: brand new solved entry
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "Sliding Window"
$ws.Cells.Item(58, 3).Value = "Fruits Into Basket"
$ws.Cells.Item(58, 4).Value = "Medium"
$ws.Cells.Item(58, 5).Value = "Done"
$ws.Cells.Item(58, 6).Value = 45898
$ws.Range("F55").Copy() | Out-Null
$ws.Range("F58").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(58, 7).Value = "O(n)"
$ws.Cells.Item(58, 8).Value = "O(1)"
$ws.Cells.Item(58, 9).Value = "Sliding Window"

# Rows 59 & 60: blank spacer rows, only the ID column is filled in.
# Clear out the stale data these rows used to hold before the insert.
$ws.Range("B59:J60").Clear() | Out-Null
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(60, 1).Value = 59

# Row 61: the still-pending "First Negative Integer..." entry, moved down
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "Sliding Window"
$ws.Cells.Item(61, 3).Value = "First Negative Integer in Every Window of Size K"
$ws.Cells.Item(61, 4).Value = "Medium"
$ws.Cells.Item(61, 5).Value = "To Do"

# Row 62: the still-pending "Count Occurrences of Anagrams" entry, moved down
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "Sliding Window"
$ws.Cells.Item(62, 3).Value = "Count Occurrences of Anagrams"
$ws.Cells.Item(62, 4).Value = "Medium"
$ws.Cells.Item(62, 5).Value = "To Do"

# Row 63: the still-pending "Majority Element" entry, moved down
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "Moore’s Voting"
$ws.Cells.Item(63, 3).Value = "Majority Element"
$ws.Cells.Item(63, 4).Value = "Easy"
$ws.Cells.Item(63, 5).Value = "To Do"

# Row 64: the still-pending "Majority Element II" entry, moved down
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "Moore’s Voting"
$ws.Cells.Item(64, 3).Value = "Majority Element II"
$ws.Cells.Item(64, 4).Value = "Medium"
$ws.Cells.Item(64, 5).Value = "To Do"

# Row 65: the still-pending "Segregate 0s, 1s, and 2s" entry, moved down
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "DNF Algorithm"
$ws.Cells.Item(65, 3).Value = "Segregate 0s, 1s, and 2s"
$ws.Cells.Item(65, 4).Value = "Easy"
$ws.Cells.Item(65, 5).Value = "To Do"

# Row 66: new trailing blank spacer row, only the ID column is filled in
$ws.Cells.Item(66, 1).Value = 65

# Update the active selection to match where the user was last editing
$ws.Range("E58:I58").Select() | Out-Null
